$d = $word.ActiveDocument

# The document currently has every run (and every paragraph mark) highlighted
# in yellow. The requested change removes all of that yellow highlighting
# throughout the whole document, leaving the text formatting otherwise
# untouched.
#
# wdNoHighlight = 0
$wdNoHighlight = 0

foreach ($p in $d.Paragraphs) {
    # Clears the highlight on every run of text within the paragraph.
    $p.Range.HighlightColorIndex = $wdNoHighlight
    # Clears the highlight carried by the paragraph mark itself (this is
    # stored separately from the runs' own formatting).
    $p.Font.HighlightColorIndex = $wdNoHighlight
}

Write-Output "Removed yellow highlighting from the document."
